$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the sample name (A column) and first scandir entry (row 2)
$ws.Range("A2").Value = "cerebellum_8rings"
$ws.Range("B2").Value = "064_GMB_ringscan_500mA_Cu100_gap5p7_orca_10x_dist50mm_100ms"

# Fill in rows 3-9 with the 8-ring data
$scandirs = @(
  "064_GMB_ringscan_500mA_Cu100_gap5p7_orca_10x_dist50mm_100ms",
  "065_GMB_ringscan_500mA_Cu100_gap5p7_orca_10x_dist50mm_100ms",
  "066_GMB_ringscan_500mA_Cu100_gap5p7_orca_10x_dist50mm_100ms",
  "067_GMB_ringscan_500mA_Cu100_gap5p7_orca_10x_dist50mm_100ms",
  "068_GMB_ringscan_500mA_Cu100_gap5p7_orca_10x_dist50mm_100ms",
  "069_GMB_ringscan_500mA_Cu100_gap5p7_orca_10x_dist50mm_100ms",
  "070_GMB_ringscan_500mA_Cu100_gap5p7_orca_10x_dist50mm_100ms"
)

for ($i = 0; $i -lt $scandirs.Length; $i++) {
  $row = 3 + $i
  $ws.Range("A$row").Value = "cerebellum_8rings"
  $ws.Range("B$row").Value = $scandirs[$i]
  $ws.Range("C$row").NumberFormat = "General"
  $ws.Range("C$row").Value = 1
  $ws.Range("D$row").Value = $i + 1
  $ws.Range("E$row").Value = $i + 2
}

# Update selection to match the authored state
$ws.Range("B12").Select()

$wb.Save()
